# items_list.xlsx edit: insert a new "Item_Name" (English name) column
# before the existing data, rename the old "Item" header to "Items",
# and refresh the AutoFilter / _FilterDatabase defined name to the new
# B1:D1 header range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank column at A - this shifts the existing
#    Item/Weights/Price columns from A:C to B:D (values, styles and
#    the autoFilter element all shift with it).
$ws.Columns("A").Insert()

# 2) Header row.
$ws.Range("A1").Value = "Item_Name"
$ws.Range("B1").Value = "Items"

# 3) English item names for column A, rows 2-22 (paired with the
#    existing Kannada names that are now in column B).
$englishNames = @(
    "Jeera",
    "Musturd",
    "Ground nut",
    "Fried gram",
    "Ragi",
    "green peas",
    "white peas",
    "Agasi seeds",
    "Ajwain",
    "Alvi",
    "Rice Flour",
    "Ragi Flour",
    "Jawar Flour",
    "Antu",
    "Badam",
    "Baking soda",
    "Chana",
    "Chana dal",
    "chana Flour",
    "Turmuric powder",
    "Pepper"
)

$row = 2
foreach ($name in $englishNames) {
    $ws.Cells.Item($row, 1).Value = $name
    $row = $row + 1
}

# 4) Column width for the new name column A (the B/C/D columns keep
#    their original widths automatically since Columns("A").Insert()
#    shifts them without touching their stored width). The COM layer
#    quantizes ColumnWidth to whole pixels, so 13.9 is the closest
#    input that lands on the saved file's target width for column A.
$ws.Columns("A").ColumnWidth = 13.9

# 5) Re-point the AutoFilter at the new header range (B:D no longer
#    includes the plain name column A).
$ws.AutoFilterMode = $false
$ws.Range("B1:D1").AutoFilter()

# 6) The _xlnm._FilterDatabase defined name still points at the old
#    range after the column insert - update it explicitly.
$filterName = $wb.Names.Item("_xlnm._FilterDatabase")
$filterName.RefersTo = "=Sheet1!`$B`$1:`$D`$1"

# 7) Restore the selection to match the saved view state.
$ws.Range("B17").Select()
